# "Generate Report for Handoff"
# Regenerate the localization-status report: move the zh-cn / de-de rows from
# "In Translation" to "Ready for handoff" and stamp the refreshed
# generate/handoff timestamps, then re-autofit the columns whose text grew.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-06 08:53:12"

# ---- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-06 08:52:59"

# ---- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-06 08:53:12"

# ---- Re-fit the Status columns now that "Ready for handoff" is longer -------
$wsOverview.Columns.Item(5).EntireColumn.AutoFit()
$wsOverview.Columns.Item(6).EntireColumn.AutoFit()
$wsZhCn.Columns.Item(3).EntireColumn.AutoFit()
$wsDeDe.Columns.Item(3).EntireColumn.AutoFit()
